# Applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.629.62"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "2.269.14"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'120.44"
$ws.Range("E5").Value = "  +6.87%  "
$ws.Range("D6").Value = "'268.84"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +2.93%  "
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").Value = "'47.43"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'0.0941"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("E12").Value = "  +5.36%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'15.76"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").Value = "2.609.65"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "2.270.51"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "43.595.50"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'6.92"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "'72.69"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("E22").Value = "  -4.67%  "
$ws.Range("D23").Value = "'234.10"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").Value = "'9.70"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "'12.30"
$ws.Range("E26").Value = "  +8.44%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").Value = "'42.03"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "'175.05"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").Value = "'21.54"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  +13.61%  "
$ws.Range("D36").Value = "'0.131"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("E37").Value = "  +7.36%  "
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.243"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'13.79"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "'72.58"
$ws.Range("E43").Value = "  -6.84%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("E46").Value = "  -9.32%  "
$ws.Range("D47").Value = "'77.23"
$ws.Range("E47").Value = "  +38.91%  "
$ws.Range("D48").Value = "'0.671"
$ws.Range("E48").Value = "  +19.37%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.58"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.27"
$ws.Range("E50").Value = "  +1.55%  "
